# dlgProfil.xlsx - "Check defaults" row gains explicit selected/not-selected
# keywords for the Profil/Abwesenheiten/Benachrichtigungen page tabs, and the
# xpath selector strings in the header row are corrected (stray `')]` -> `']`,
# and the double-quoted Benachrichtigungen variant normalized to match the
# others).
#
# NOTE: shared-string indices in the saved file are allocated in the order
# cells are first written, so row 3 (<SELECTED>/<NOTSELECTED>) is populated
# before row 2 (the xpath strings) to reproduce the same string order as the
# target workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 ("Check defaults"): Profil tab is selected by default, the other two
# are not.
$ws.Range("C3").Value = "<SELECTED>"
$ws.Range("D3").Value = "<NOTSELECTED>"
$ws.Range("E3").Value = "<NOTSELECTED>"

# Row 2 ("Record/Selector"): fix the xpath selectors for the three nav-item
# links.
$ws.Range("C2").Value = "xpath=//*/li[@class = 'nav-item']/a[text()=' Profil']"
$ws.Range("D2").Value = "xpath=//*/li[@class = 'nav-item']/a[text()=' Abwesenheiten']"
$ws.Range("E2").Value = "xpath=//*/li[@class = 'nav-item']/a[text()=' Benachrichtigungen']"

# Page setup: paper size 9 = A4, orientation 1 = portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Leave the cursor on E3, matching the saved selection.
$ws.Range("E3").Select()
